# Cotações atualizadas - 2025-10-20
# Add a new data row (row 46) to the quotes sheet with the latest values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 46

# Date serial 45950 == 2025-10-20
$ws.Cells.Item($newRow, 1).Value = 45950
$ws.Cells.Item($newRow, 2).Value = "21,4368"
$ws.Cells.Item($newRow, 3).Value = "15,3478"
$ws.Cells.Item($newRow, 4).Value = "15,4273"
$ws.Cells.Item($newRow, 5).Value = "15,4273"

# Match the date formatting/style used by the previous row's date cell (A45).
$ws.Range("A45").Copy()
$ws.Range("A46").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

Write-Output "Row 46 added: 2025-10-20 quotes"
